$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.944.37'
$ws.Range("E2").Value = '  +4.22%  '
$ws.Range("D3").Value = '2.235.10'
$ws.Range("E3").Value = '  +2.77%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.49'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0931'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.50'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").Value = '2.551.31'
$ws.Range("E15").Value = '  +2.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").Value = '2.235.07'
$ws.Range("E17").Value = '  +3.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.810'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = '42.870.09'
$ws.Range("E19").Value = '  +4.49%  '
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '230.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.20'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +12.27%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -5.44%  '
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.75'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.82%  '
$ws.Range("E32").Value = '  +20.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0792'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.23%  '
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  +7.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0329'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +15.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.13'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.56'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '59.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '105.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0989'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.53%  '
$ws.Range("E48").Value = '  +19.48%  '
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.13%  '
